$wb = $excel.ActiveWorkbook

# Map of 1-based worksheet index -> new cell values to write into column A.
# Sheet 1 (Compartment) holds BOTH the document-level '!!!ObjTables ...' banner (A1)
# and the table-level '!!ObjTables ...' banner (A2); every other sheet holds only
# its own table-level banner in A1.
$sheetEdits = @(
    @{ Index = 1; A1 = '!!!ObjTables objTablesVersion=''0.0.8'' date=''2020-03-09 13:01:27'''; A2 = '!!ObjTables type=''Data'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8''' }
    @{ Index = 2; A1 = '!!ObjTables type=''Data'' id=''Compound'' name=''Compound'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 3; A1 = '!!ObjTables type=''Data'' id=''Definition'' name=''Definition'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 4; A1 = '!!ObjTables type=''Data'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 5; A1 = '!!ObjTables type=''Data'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 6; A1 = '!!ObjTables type=''Data'' id=''Gene'' name=''Gene'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 7; A1 = '!!ObjTables type=''Data'' id=''Layout'' name=''Layout'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 8; A1 = '!!ObjTables type=''Data'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 9; A1 = '!!ObjTables type=''Data'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 10; A1 = '!!ObjTables type=''Data'' id=''Position'' name=''Position'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 11; A1 = '!!ObjTables type=''Data'' id=''Protein'' name=''Protein'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 12; A1 = '!!ObjTables type=''Data'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 13:01:27'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 13; A1 = '!!ObjTables type=''Data'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 14; A1 = '!!ObjTables type=''Data'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 15; A1 = '!!ObjTables type=''Data'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 16; A1 = '!!ObjTables type=''Data'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 17; A1 = '!!ObjTables type=''Data'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 18; A1 = '!!ObjTables type=''Data'' id=''Relation'' name=''Relation'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 19; A1 = '!!ObjTables type=''Data'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 20; A1 = '!!ObjTables type=''Data'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 21; A1 = '!!ObjTables type=''Data'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 22; A1 = '!!ObjTables type=''Data'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 23; A1 = '!!ObjTables type=''Data'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 24; A1 = '!!ObjTables type=''Data'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 25; A1 = '!!ObjTables type=''Data'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
    @{ Index = 26; A1 = '!!ObjTables type=''Data'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 13:01:28'' objTablesVersion=''0.0.8'''; A2 = $null }
)

foreach ($edit in $sheetEdits) {
    $ws = $wb.Worksheets.Item($edit.Index)
    # The sheets ship protected (no password); unprotect, write, then restore
    # protection so the sheet ends up locked again, matching the source document.
    $ws.Unprotect()
    $ws.Range("A1").Value = $edit.A1
    if ($edit.A2) {
        $ws.Range("A2").Value = $edit.A2
    }
    $ws.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $false, $false, $false)
}
